$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the description text for the third log entry (row 4, column E)
$ws.Range("E4").Value = "Finished Chapter 5. Learned more about debugging. This is something I was not as familiar with as I thought so it was a good thing I learned some proper practices, and Python specific ones too. Also finished Chapter 6. This was all about lists. It had the first proper programming challenge, writing a pragram called coinFlipStreak.py. From now on I expect to code more and learn that way. Forgot to upload some files so updated next day."

# Update the "Time end" value for that same row from 17:00 to 18:00
$ws.Range("C4").Value = 0.75
